$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the province names between rows 53 and 54 (Huelva/Huesca reorder
# in the shared-string table manifests here as the two rows trading labels)
$ws.Range("A53").Value = "Huesca"
$ws.Range("A54").Value = "Huelva"

# Swap "Casos activos" values to match the reordering
$ws.Range("C53").Value = 0
$ws.Range("C54").Value = 72

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 06:46"
